# Fix the typo in the "Mark as issue" requirement description and restore
# the last-used cell selection, matching the author's "Added readme and
# updated docs" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Requirement Description for the "Mark as issue" row (D2): "use" -> "user"
$ws.Range("D2").Value = "Any user logged in can mark an issue in the project."

# Active cell moved from E6 to A6 before the file was saved.
$ws.Range("A6").Select() | Out-Null
